# Insert a new weekly record for "Femacal de La Calera" / Cebollín at row 762.
# Excel shifts the previous rows 762-799 down to 763-800 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(762).Insert()

$ws.Range("A762").Value = 3
$ws.Range("B762").Value = "Femacal de La Calera"
$ws.Range("C762").Value = "Coquimbo"
$ws.Range("D762").Value = 45147
$ws.Range("E762").Value = 5
$ws.Range("F762").Value = 100112037
$ws.Range("G762").Value = "Cebollín"
$ws.Range("H762").Value = "Sin especificar"
$ws.Range("I762").Value = "Primera"
$ws.Range("J762").Value = 120
$ws.Range("K762").Value = 4000
$ws.Range("L762").Value = 4000
$ws.Range("M762").Value = 4000
$ws.Range("N762").Value = "`$/paquete 36 unidades"
$ws.Range("O762").Value = "Provincia de Quillota"
$ws.Range("P762").Value = 111
$ws.Range("Q762").Value = 36
$ws.Range("R762").Value = "Hortaliza"
